# Updates the cryptos list with refreshed price / volume(1h) figures,
# and fixes the swapped NEARProtocol / ARBITRUM rows (43 & 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the Price (D) and Volume(1h) (E) columns keep being plain text
# (some of the new values, e.g. "399.63", look numeric and Excel would
# otherwise silently convert them to numbers instead of leaving them as
# the literal strings found in the source data).
$ws.Range("B2:E51").NumberFormat = "@"

$changes = @{
    'D2' = '58.592.94'
    'E2' = '  +4.05%  '
    'D3' = '3.297.28'
    'E3' = '  +2.16%  '
    'E4' = '  -0.02%  '
    'D5' = '399.63'
    'E5' = '  +0.44%  '
    'D6' = '109.61'
    'E6' = '  -1.48%  '
    'E7' = '  +5.46%  '
    'E8' = '  +0.05%  '
    'D9' = '0.634'
    'E9' = '  +2.48%  '
    'D10' = '39.78'
    'E10' = '  +1.19%  '
    'D11' = '0.0983'
    'E11' = '  +7.44%  '
    'E12' = '  +1.43%  '
    'D13' = '3.821.79'
    'E13' = '  +2.16%  '
    'D14' = '8.35'
    'E14' = '  +3.01%  '
    'D15' = '19.24'
    'E15' = '  +0.94%  '
    'D16' = '3.299.18'
    'E16' = '  +2.48%  '
    'E17' = '  -0.44%  '
    'D18' = '10.90'
    'E18' = '  -0.14%  '
    'D19' = '58.333.29'
    'E19' = '  +3.79%  '
    'D20' = '3.31'
    'E20' = '  -0.99%  '
    'D21' = '0.0000110'
    'E21' = '  +6.77%  '
    'D22' = '12.94'
    'E22' = '  -0.59%  '
    'D23' = '302.03'
    'E23' = '  +1.16%  '
    'D24' = '74.58'
    'E24' = '  -1.60%  '
    'D25' = '3.18'
    'E25' = '  -0.94%  '
    'D26' = '28.33'
    'E26' = '  +0.83%  '
    'D27' = '4.42'
    'E27' = '  +1.06%  '
    'D28' = '7.87'
    'E28' = '  -3.41%  '
    'D29' = '7.40'
    'E29' = '  -0.62%  '
    'D30' = '0.171'
    'E30' = '  -1.49%  '
    'E31' = '  -0.46%  '
    'E32' = '  +2.71%  '
    'D33' = '11.37'
    'E33' = '  +2.07%  '
    'E34' = '  +11.04%  '
    'D35' = '0.0530'
    'E35' = '  +7.41%  '
    'E36' = '  +0.28%  '
    'D37' = '51.86'
    'E37' = '  +0.94%  '
    'D38' = '3.28'
    'E38' = '  +4.81%  '
    'E39' = '  -0.06%  '
    'D40' = '3.48'
    'E40' = '  -1.31%  '
    'D41' = '137.80'
    'E41' = '  +0.48%  '
    'E42' = '  +2.71%  '
    'B43' = 'ARBITRUM'
    'C43' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D43' = '1.88'
    'E43' = '  -1.98%  '
    'B44' = 'NEARProtocol'
    'C44' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D44' = '3.92'
    'E44' = '  -2.26%  '
    'D45' = '16.78'
    'E45' = '  -3.90%  '
    'E46' = '  -1.94%  '
    'D47' = '2.29'
    'E47' = '  +9.63%  '
    'D48' = '22.46'
    'E48' = '  +0.81%  '
    'D49' = '2.161.87'
    'E49' = '  +1.44%  '
    'E50' = '  -0.50%  '
    'E51' = '  -13.66%  '
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}
